# Update contrib_before_closest_residue (P), contrib_closest_residue (Q),
# and contrib_next_closest_residue (R) columns with recomputed residue
# distances based on charges for the 6799 frames results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 130
$ws.Range("Q2").Value = 130
$ws.Range("R2").Value = 130
$ws.Range("P4").Value = 1105
$ws.Range("Q4").Value = 1105
$ws.Range("R4").Value = 1105
$ws.Range("P5").Value = 780
$ws.Range("Q5").Value = 780
$ws.Range("Q8").Value = 455
$ws.Range("R8").Value = 455
$ws.Range("Q9").Value = 780
$ws.Range("Q10").Value = 130
$ws.Range("R10").Value = 130
$ws.Range("Q16").Value = 130
$ws.Range("R16").Value = 130
$ws.Range("P17").Value = 780
$ws.Range("Q17").Value = 780
$ws.Range("R17").Value = 780
$ws.Range("P22").Value = 1105
$ws.Range("P23").Value = 130
$ws.Range("R23").Value = 130
$ws.Range("P25").Value = 1105
$ws.Range("Q25").Value = 1105
$ws.Range("R25").Value = 1105
$ws.Range("P26").Value = 780
$ws.Range("Q26").Value = 780
$ws.Range("R26").Value = 780
$ws.Range("P27").Value = 455
$ws.Range("Q27").Value = 455
$ws.Range("R27").Value = 455
$ws.Range("P28").Value = 130
$ws.Range("Q28").Value = 130
$ws.Range("R28").Value = 130
$ws.Range("P29").Value = 1105
$ws.Range("Q29").Value = 1105
$ws.Range("R29").Value = 1105
$ws.Range("Q31").Value = "SF"
$ws.Range("R31").Value = 1105
$ws.Range("Q32").Value = 780
$ws.Range("P33").Value = 455
$ws.Range("Q33").Value = 455
$ws.Range("R33").Value = 455
$ws.Range("P34").Value = 1105
$ws.Range("Q34").Value = 1105
$ws.Range("R34").Value = 1105
$ws.Range("P37").Value = 1105
$ws.Range("Q37").Value = 1105
$ws.Range("R37").Value = 1105
$ws.Range("Q38").Value = 780
$ws.Range("R38").Value = 780
$ws.Range("P39").Value = 1105
$ws.Range("Q39").Value = 1105
$ws.Range("R39").Value = 1105
$ws.Range("P40").Value = 780
$ws.Range("Q40").Value = 780
$ws.Range("R40").Value = 780
$ws.Range("P41").Value = 130
$ws.Range("Q41").Value = 130
$ws.Range("P45").Value = 1105
$ws.Range("Q45").Value = 1105
$ws.Range("R45").Value = 1105
$ws.Range("P46").Value = 780
$ws.Range("Q46").Value = 780
$ws.Range("R46").Value = 780
$ws.Range("P47").Value = 455
$ws.Range("Q47").Value = 455
$ws.Range("R47").Value = 455
$ws.Range("P49").Value = 1105
$ws.Range("Q49").Value = 1105
$ws.Range("R49").Value = 1105
$ws.Range("P50").Value = 455
$ws.Range("Q50").Value = 455
$ws.Range("R50").Value = 455
$ws.Range("P51").Value = 780
$ws.Range("Q51").Value = 780
$ws.Range("R51").Value = 780
$ws.Range("P52").Value = 455
$ws.Range("Q52").Value = 455
$ws.Range("R52").Value = 455
$ws.Range("P55").Value = 1105
$ws.Range("Q55").Value = 1105
$ws.Range("R55").Value = 1105
$ws.Range("P56").Value = 780
$ws.Range("Q56").Value = 780
$ws.Range("R56").Value = 780
$ws.Range("P57").Value = 1105
$ws.Range("Q57").Value = 1105
$ws.Range("P59").Value = 780
$ws.Range("Q59").Value = 780
$ws.Range("R59").Value = 780
$ws.Range("P60").Value = 130
$ws.Range("Q60").Value = 455
$ws.Range("R60").Value = 455
$ws.Range("R62").Value = 455
$ws.Range("Q63").Value = 130
$ws.Range("P64").Value = 455
$ws.Range("Q64").Value = 130
$ws.Range("R64").Value = 455
$ws.Range("P66").Value = 1105
$ws.Range("Q66").Value = 1105
$ws.Range("R66").Value = 1105
$ws.Range("P67").Value = 780
$ws.Range("Q67").Value = 780
$ws.Range("R67").Value = 780
$ws.Range("Q71").Value = 780
$ws.Range("R71").Value = 780
$ws.Range("Q72").Value = 1105
$ws.Range("R72").Value = 1105
$ws.Range("P73").Value = 455
$ws.Range("Q73").Value = 455
$ws.Range("R73").Value = 455
$ws.Range("P74").Value = 455
$ws.Range("Q74").Value = 455
$ws.Range("R74").Value = 455
$ws.Range("Q75").Value = 130
$ws.Range("R75").Value = 130
$ws.Range("P77").Value = 130
$ws.Range("Q77").Value = 130
$ws.Range("R77").Value = 130
$ws.Range("P78").Value = 1105
$ws.Range("Q78").Value = 1105
$ws.Range("R78").Value = 1105
$ws.Range("R79").Value = 455
$ws.Range("P80").Value = 1105
$ws.Range("Q80").Value = 1105
$ws.Range("P82").Value = 455
$ws.Range("Q82").Value = 455
$ws.Range("R82").Value = 455
$ws.Range("P84").Value = 1105
$ws.Range("Q84").Value = 1105
$ws.Range("R84").Value = 1105
$ws.Range("P86").Value = 130
$ws.Range("P88").Value = 130
$ws.Range("Q88").Value = 130
$ws.Range("R88").Value = 130
$ws.Range("Q89").Value = 130
$ws.Range("P90").Value = 455
$ws.Range("P91").Value = 1105
$ws.Range("P92").Value = 130
$ws.Range("Q92").Value = 130
$ws.Range("R92").Value = 130
$ws.Range("P94").Value = 455
$ws.Range("Q94").Value = 455
$ws.Range("R94").Value = 455
$ws.Range("Q95").Value = 1105
$ws.Range("R95").Value = 1105
$ws.Range("P97").Value = 455
$ws.Range("P98").Value = "SF"
$ws.Range("Q98").Value = "SF"
$ws.Range("R98").Value = 130
$ws.Range("P99").Value = 130
$ws.Range("Q99").Value = 130
$ws.Range("R99").Value = 130
